$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: EnCQR-LSTM -> Proporcion_Sig "4/10" becomes "3/10", Mejor_N_Calib (C2) 102.4 -> 76.8
$ws.Range("B2").Value = "3/10"
$ws.Range("C2").Value = 76.8

# Row 3: AREPD -> Proporcion_Sig "1/10" becomes "0/10", Mejor_N_Calib (C3) 25.6 -> 0
$ws.Range("B3").Value = "0/10"
$ws.Range("C3").Value = 0
